$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing rows that match the desired style:
#  - Column A (data structure tag) uses the plain bordered/wrapped style (like A29)
#  - Column B (question title) uses the yellow-filled bordered/wrapped style (like B28)
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B28").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in the new "pending" question row
$ws.Cells.Item(30, 1).Value2 = "Linked List"
$ws.Cells.Item(30, 2).Value2 = "143. Reorder List"

# Update the selection to match the post-edit cursor position
$ws.Range("C35").Select() | Out-Null
